$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are stored as text in the source workbook (t="inlineStr"),
# e.g. "217.29" or "1.010". Excel would normally auto-convert such strings to numbers,
# so we force the NumberFormat to text ("@") before assigning these values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.202.33"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.659.17"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.29"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5172"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2640"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06274"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.81"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07771"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.489"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.667.52"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.886.01"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5471"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8134"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.90"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.209.05"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.610"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.39"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.09"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.002"
$ws.Range("E23").Value = "  -4.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.45"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1220"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.284"
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.20"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.440"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05936"
$ws.Range("E30").Value = "  -4.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.275"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.549"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.285"
$ws.Range("E33").Value = "  -4.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.584"
$ws.Range("E34").Value = "  -6.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9608"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.420"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.768"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5688"
$ws.Range("E38").Value = "  -6.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.042"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01592"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8543"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.014.41"
$ws.Range("E43").Value = "  -6.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.04"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.799.85"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.54"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.010"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.040"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05168"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4226"
$ws.Range("E51").Value = "  -0.04%  "
